$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated hours worked on Friday (column F) for all 5 team members in week 1
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 2
$ws.Range("F7").Value = 2

# Move active selection to L7 (side effect of editing session)
$ws.Range("L7").Select()
